# Apply text replacements to the two-digit multiplication table.
$d = $word.ActiveDocument

$replacements = @(
    @{old = "90×95=8550"; new = "44×91=4004"},
    @{old = "24×22=528";  new = "44×42=1848"},
    @{old = "72×94=6768"; new = "90×81=7290"},
    @{old = "58×14=812";  new = "86×43=3698"},
    @{old = "62×48=2976"; new = "83×30=2490"},
    @{old = "13×20=260";  new = "67×73=4891"},
    @{old = "77×71=5467"; new = "70×93=6510"},
    @{old = "91×45=4095"; new = "69×40=2760"},
    @{old = "97×58=5626"; new = "55×76=4180"},
    @{old = "92×99=9108"; new = "16×47=752"},
    @{old = "51×86=4386"; new = "53×23=1219"},
    @{old = "29×59=1711"; new = "98×52=5096"},
    @{old = "14×63=882";  new = "56×90=5040"},
    @{old = "89×58=5162"; new = "71×20=1420"},
    @{old = "32×21=672";  new = "93×23=2139"},
    @{old = "96×49=4704"; new = "89×76=6764"},
    @{old = "62×57=3534"; new = "21×45=945"},
    @{old = "42×31=1302"; new = "92×11=1012"},
    @{old = "19×29=551";  new = "86×64=5504"},
    @{old = "49×81=3969"; new = "74×65=4810"},
    @{old = "25×65=1625"; new = "48×26=1248"},
    @{old = "47×12=564";  new = "73×55=4015"},
    @{old = "43×32=1376"; new = "92×68=6256"},
    @{old = "11×60=660";  new = "89×79=7031"},
    @{old = "79×36=2844"; new = "85×80=6800"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
